$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7331460674157303
$ws1.Range("C2").Value = 0.8738738738738738
$ws1.Range("D2").Value = 0.5449438202247191
$ws1.Range("E2").Value = 0.671280276816609
$ws1.Range("F2").Value = 0.5893074119076549
$ws1.Range("G2").Value = 0.5529489147116861
$ws1.Range("H2").Value = 0.7331460674157304
$ws1.Range("I2").Value = 291
$ws1.Range("J2").Value = 42
$ws1.Range("K2").Value = 492
$ws1.Range("L2").Value = 243

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.6693877551020408
$ws2.Range("C2").Value = 0.9213483146067416
$ws2.Range("D2").Value = 0.7754137115839244

$ws2.Range("B3").Value = 0.8738738738738738
$ws2.Range("C3").Value = 0.5449438202247191
$ws2.Range("D3").Value = 0.671280276816609

$ws2.Range("B4").Value = 0.7331460674157303
$ws2.Range("C4").Value = 0.7331460674157303
$ws2.Range("D4").Value = 0.7331460674157303
$ws2.Range("E4").Value = 0.7331460674157303

$ws2.Range("B5").Value = 0.7716308144879573
$ws2.Range("C5").Value = 0.7331460674157304
$ws2.Range("D5").Value = 0.7233469942002666

$ws2.Range("B6").Value = 0.7716308144879573
$ws2.Range("C6").Value = 0.7331460674157303
$ws2.Range("D6").Value = 0.7233469942002667

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 492
$ws3.Range("C2").Value = 42
$ws3.Range("B3").Value = 243
$ws3.Range("C3").Value = 291
